$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1033
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13
$ws.Range("H18").Value = 111111310
$ws.Range("I18").Value = 111111310
$ws.Range("K18").Value = 111111310
$ws.Range("M18").Value = -111111026
$ws.Range("H28").Value = 1898
$ws.Range("I28").Value = 837
$ws.Range("J28").Value = 4444.4
$ws.Range("K28").Value = 837
$ws.Range("L28").Value = 4444.4
$ws.Range("M28").Value = -352
$ws.Range("N28").Value = -5414.4
$ws.Range("H40").Value = 5578280.5
$ws.Range("J40").Value = 11134610
$ws.Range("L40").Value = 11134610
$ws.Range("N40").Value = -11134960
$ws.Range("H43").Value = 1934951.9
$ws.Range("I43").Value = 3856903.8
$ws.Range("K43").Value = 3856903.8
$ws.Range("M43").Value = -3856834.8
$ws.Range("H64").Value = 23633866
$ws.Range("J64").Value = 4894.2
$ws.Range("L64").Value = 4894.2
$ws.Range("N64").Value = -5390.2
$ws.Range("H67").Value = 23633866
$ws.Range("J67").Value = 4894.2
$ws.Range("L67").Value = 4894.2
$ws.Range("N67").Value = -6610.2
$ws.Range("H112").Value = 4153.7144
$ws.Range("J112").Value = 4153.7144
$ws.Range("L112").Value = 12461.1432
$ws.Range("N112").Value = -14677.1432
$ws.Range("H127").Value = 2426.186
$ws.Range("I127").Value = 829.6667
$ws.Range("J127").Value = 3281.4644
$ws.Range("K127").Value = 2489.0001
$ws.Range("L127").Value = 9844.393199999999
$ws.Range("M127").Value = 2470.9999
$ws.Range("N127").Value = -19764.3932
$ws.Range("H129").Value = 3545.2354
$ws.Range("J129").Value = 6137.8887
$ws.Range("L129").Value = 18413.6661
$ws.Range("N129").Value = -28413.6661
$ws.Range("H132").Value = 222784.83
$ws.Range("I132").Value = 258777.75
$ws.Range("K132").Value = 776333.25
$ws.Range("M132").Value = -773803.25
$ws.Range("H137").Value = 5498.0557
$ws.Range("I137").Value = 4199.143
$ws.Range("K137").Value = 12597.429
$ws.Range("M137").Value = -10047.429
$ws.Range("H138").Value = 7698.797
$ws.Range("J138").Value = 8576.925999999999
$ws.Range("L138").Value = 25730.778
$ws.Range("N138").Value = -36010.778
$ws.Range("H141").Value = 6753.8184
$ws.Range("I141").Value = 7411.5
$ws.Range("K141").Value = 22234.5
$ws.Range("M141").Value = -17054.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2114.5583
$ws.Range("I32").Value = 1689.5286
$ws.Range("K32").Value = 1689.5286
$ws.Range("M32").Value = -1402.5286
$ws.Range("H61").Value = 11246.257
$ws.Range("J61").Value = 14025
$ws.Range("L61").Value = 14025
$ws.Range("N61").Value = -14449
$ws.Range("H74").Value = 3841.574
$ws.Range("I74").Value = 908.25
$ws.Range("J74").Value = 5076.6577
$ws.Range("K74").Value = 908.25
$ws.Range("L74").Value = 5076.6577
$ws.Range("M74").Value = -34.25
$ws.Range("N74").Value = -6824.6577
$ws.Range("H77").Value = 3841.574
$ws.Range("I77").Value = 908.25
$ws.Range("J77").Value = 5076.6577
$ws.Range("K77").Value = 4541.25
$ws.Range("L77").Value = 25383.2885
$ws.Range("M77").Value = -173.25
$ws.Range("N77").Value = -34119.2885
$ws.Range("H136").Value = 11246.257
$ws.Range("J136").Value = 14025
$ws.Range("L136").Value = 42075
$ws.Range("N136").Value = -47175

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2931.25
$ws.Range("I86").Value = 2866
$ws.Range("J86").Value = 3127
$ws.Range("K86").Value = 2866
$ws.Range("L86").Value = 3127
$ws.Range("M86").Value = -1743
$ws.Range("N86").Value = -5373
$ws.Range("H89").Value = 2931.25
$ws.Range("I89").Value = 2866
$ws.Range("J89").Value = 3127
$ws.Range("K89").Value = 14330
$ws.Range("L89").Value = 15635
$ws.Range("M89").Value = -8714
$ws.Range("N89").Value = -26867
$ws.Range("H105").Value = 1193.8
$ws.Range("I105").Value = 996.3333
$ws.Range("K105").Value = 996.3333
$ws.Range("M105").Value = 750.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2526474
$ws.Range("I105").Value = 5682617
$ws.Range("K105").Value = 5682617
$ws.Range("M105").Value = -5680870

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 237588.84
$ws.Range("J68").Value = 279694.03
$ws.Range("L68").Value = 839082.0900000001
$ws.Range("N68").Value = -840704.0900000001
$ws.Range("H71").Value = 237588.84
$ws.Range("J71").Value = 279694.03
$ws.Range("L71").Value = 2517246.27
$ws.Range("N71").Value = -2525358.27
$ws.Range("H112").Value = 12747.5
$ws.Range("I112").Value = 495
$ws.Range("K112").Value = 1485
$ws.Range("M112").Value = -377
$ws.Range("H131").Value = 1546077.5
$ws.Range("J131").Value = 1982441.4
$ws.Range("L131").Value = 5947324.199999999
$ws.Range("N131").Value = -5957404.199999999
$ws.Range("H132").Value = 1879.1333
$ws.Range("I132").Value = 468.7
$ws.Range("J132").Value = 4700
$ws.Range("K132").Value = 4218.3
$ws.Range("L132").Value = 42300
$ws.Range("M132").Value = -1688.3
$ws.Range("N132").Value = -47360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 22390.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 62500228
$ws.Range("I55").Value = 83333490
$ws.Range("J55").Value = 450
$ws.Range("K55").Value = 83333490
$ws.Range("L55").Value = 450
$ws.Range("M55").Value = -83333317
$ws.Range("N55").Value = -796
$ws.Range("H61").Value = 3642.1177
$ws.Range("I61").Value = 3570.2068
$ws.Range("J61").Value = 4059.2
$ws.Range("K61").Value = 3570.2068
$ws.Range("L61").Value = 4059.2
$ws.Range("M61").Value = -3368.2068
$ws.Range("N61").Value = -4463.2
$ws.Range("H113").Value = 3642.1177
$ws.Range("I113").Value = 3570.2068
$ws.Range("J113").Value = 4059.2
$ws.Range("K113").Value = 3570.2068
$ws.Range("L113").Value = 4059.2
$ws.Range("M113").Value = -1400.2068
$ws.Range("N113").Value = -8399.200000000001
$ws.Range("H136").Value = 4204.05
$ws.Range("I136").Value = 3963.5945
$ws.Range("J136").Value = 4888.423
$ws.Range("K136").Value = 11890.7835
$ws.Range("L136").Value = 14665.269
$ws.Range("M136").Value = -9340.783500000001
$ws.Range("N136").Value = -19765.269

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 65000000
$ws.Range("I13").Value = 65000000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 65000000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -64999860
$ws.Range("N13").ClearContents()
$ws.Range("H81").Value = 1614735.2
$ws.Range("I81").Value = 1897919.5
$ws.Range("K81").Value = 3795839
$ws.Range("M81").Value = -3794778
$ws.Range("H84").Value = 1614735.2
$ws.Range("I84").Value = 1897919.5
$ws.Range("K84").Value = 18979195
$ws.Range("M84").Value = -18973891
$ws.Range("H100").Value = 1112812.8
$ws.Range("I100").Value = 1667806
$ws.Range("J100").Value = 2826.3333
$ws.Range("K100").Value = 3335612
$ws.Range("L100").Value = 5652.6666
$ws.Range("M100").Value = -3335071
$ws.Range("N100").Value = -6734.6666
$ws.Range("H122").Value = 3170.5686
$ws.Range("I122").Value = 2499.8538
$ws.Range("J122").Value = 5920.5
$ws.Range("K122").Value = 7499.5614
$ws.Range("L122").Value = 17761.5
$ws.Range("M122").Value = -5049.5614
$ws.Range("N122").Value = -22661.5
